$d = $word.ActiveDocument

# 1. Abstract paragraph: "experienced" -> "gone through"
$d.Content.Find.Execute(
    "Last twenty years, vehicle monitoring systems have experienced significant rise",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Last twenty years, vehicle monitoring systems have gone through significant rise",
    2) | Out-Null

# 2. "fleet management such as delivery trucks" -> "fleet management, in particular delivery trucks"
$d.Content.Find.Execute(
    "That kind of systems is suitable for fleet management such as delivery trucks",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "That kind of systems is suitable for fleet management, in particular delivery trucks",
    2) | Out-Null

# 3. "they use self-developed GPS tracker devices" -> "they use the self-developed GPS tracker devices"
$d.Content.Find.Execute(
    "they use self-developed GPS tracker devices",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "they use the self-developed GPS tracker devices",
    2) | Out-Null

# 4. "Event triggers, sends" -> "Event triggers sends" (remove comma)
$d.Content.Find.Execute(
    "Event triggers, sends notification",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Event triggers sends notification",
    2) | Out-Null

# 5. Merge "Service " + "Use Case diagram" runs into a single run (no visible text change)
$d.Content.Find.Execute(
    "Service Use Case diagram",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Service Use Case diagram",
    2) | Out-Null

# 6. "As suggested by this paper, while" -> "As suggested by this article, while"
$d.Content.Find.Execute(
    "As suggested by this paper, while vehicles monitoring system market is already divided by historically successful competitors, which use the same approach to ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "As suggested by this article, while vehicles monitoring system market is already divided by historically successful competitors, which use the same approach to ",
    2) | Out-Null
